$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H75").Value = 38869
$ws.Range("J75").Value = 38869
$ws.Range("L75").Value = 38869
$ws.Range("N75").Value = -40741
$ws.Range("H78").Value = 38869
$ws.Range("J78").Value = 38869
$ws.Range("L78").Value = 116607
$ws.Range("N78").Value = -125967
$ws.Range("H123").Value = 29960.445
$ws.Range("J123").Value = 29960.445
$ws.Range("L123").Value = 29960.445
$ws.Range("N123").Value = -39760.445
$ws.Range("H124").Value = 43306.2
$ws.Range("J124").Value = 43306.2
$ws.Range("L124").Value = 43306.2
$ws.Range("N124").Value = -53126.2
$ws.Range("H128").Value = 36858.5
$ws.Range("J128").Value = 36858.5
$ws.Range("L128").Value = 36858.5
$ws.Range("N128").Value = -46818.5
$ws.Range("H129").Value = 3076.7
$ws.Range("I129").Value = 7598.5
$ws.Range("J129").Value = 1946.25
$ws.Range("K129").Value = 22795.5
$ws.Range("L129").Value = 5838.75
$ws.Range("M129").Value = -17795.5
$ws.Range("N129").Value = -15838.75
$ws.Range("H130").Value = 43503.2
$ws.Range("J130").Value = 43503.2
$ws.Range("L130").Value = 43503.2
$ws.Range("N130").Value = -53543.2
$ws.Range("H132").Value = 26163.87
$ws.Range("J132").Value = 126100.43
$ws.Range("L132").Value = 378301.29
$ws.Range("N132").Value = -383361.29

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H80").Value = 49996
$ws.Range("J80").Value = 49996
$ws.Range("L80").Value = 49996
$ws.Range("N80").Value = -51992
$ws.Range("H83").Value = 49996
$ws.Range("J83").Value = 49996
$ws.Range("L83").Value = 149988
$ws.Range("N83").Value = -159972
$ws.Range("H128").Value = 49996
$ws.Range("J128").Value = 49996
$ws.Range("L128").Value = 49996
$ws.Range("N128").Value = -59956
$ws.Range("H131").Value = 49992
$ws.Range("J131").Value = 49992
$ws.Range("L131").Value = 49992
$ws.Range("N131").Value = -60072
$ws.Range("H137").Value = 49800
$ws.Range("J137").Value = 49800
$ws.Range("L137").Value = 49800
$ws.Range("N137").Value = -60000

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H117").Value = 50000
$ws.Range("J117").Value = 50000
$ws.Range("L117").Value = 50000
$ws.Range("N117").Value = -59178
$ws.Range("H122").Value = 36102.4
$ws.Range("J122").Value = 36102.4
$ws.Range("L122").Value = 36102.4
$ws.Range("N122").Value = -45902.4
$ws.Range("H126").Value = 50776
$ws.Range("J126").Value = 50776
$ws.Range("L126").Value = 50776
$ws.Range("N126").Value = -60656
$ws.Range("H130").Value = 39529
$ws.Range("J130").Value = 39529
$ws.Range("L130").Value = 39529
$ws.Range("N130").Value = -49569
$ws.Range("H138").Value = 28142.857
$ws.Range("J138").Value = 28142.857
$ws.Range("L138").Value = 28142.857
$ws.Range("N138").Value = -38422.857

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H20").Value = 45906.8
$ws.Range("J20").Value = 45906.8
$ws.Range("L20").Value = 45906.8
$ws.Range("N20").Value = -46378.8
$ws.Range("H30").Value = 45906.8
$ws.Range("J30").Value = 45906.8
$ws.Range("L30").Value = 45906.8
$ws.Range("N30").Value = -46088.8
$ws.Range("H88").Value = 44226
$ws.Range("J88").Value = 44226
$ws.Range("L88").Value = 44226
$ws.Range("N88").Value = -45038
$ws.Range("H91").Value = 44226
$ws.Range("J91").Value = 44226
$ws.Range("L91").Value = 44226
$ws.Range("N91").Value = -47034
$ws.Range("H100").Value = 37440.668
$ws.Range("J100").Value = 37440.668
$ws.Range("L100").Value = 37440.668
$ws.Range("N100").Value = -39604.668
$ws.Range("H111").Value = 46997.332
$ws.Range("J111").Value = 46997.332
$ws.Range("L111").Value = 46997.332
$ws.Range("N111").Value = -55177.332
$ws.Range("H128").Value = 45906.8
$ws.Range("J128").Value = 45906.8
$ws.Range("L128").Value = 45906.8
$ws.Range("N128").Value = -55866.8
$ws.Range("H130").Value = 38780
$ws.Range("J130").Value = 38780
$ws.Range("L130").Value = 38780
$ws.Range("N130").Value = -48820
$ws.Range("H138").Value = 50533
$ws.Range("J138").Value = 50533
$ws.Range("L138").Value = 50533
$ws.Range("N138").Value = -60813

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 5131.3706
$ws.Range("I131").Value = 25657
$ws.Range("J131").Value = 1561.6957
$ws.Range("K131").Value = 76971
$ws.Range("L131").Value = 4685.0871
$ws.Range("M131").Value = -71931
$ws.Range("N131").Value = -14765.0871
$ws.Range("H137").Value = 52641870
$ws.Range("I137").Value = 3436.25
$ws.Range("J137").Value = 90924370
$ws.Range("K137").Value = 10308.75
$ws.Range("L137").Value = 272773110
$ws.Range("M137").Value = -5208.75
$ws.Range("N137").Value = -272783310

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H130").Value = 48881.332
$ws.Range("J130").Value = 48881.332
$ws.Range("L130").Value = 48881.332
$ws.Range("N130").Value = -58921.332
$ws.Range("H139").Value = 29760
$ws.Range("J139").Value = 29760
$ws.Range("L139").Value = 29760
$ws.Range("N139").Value = -40040

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H36").Value = 48715
$ws.Range("J36").Value = 48715
$ws.Range("L36").Value = 48715
$ws.Range("N36").Value = -49839
$ws.Range("H96").Value = 30997.334
$ws.Range("J96").Value = 30997.334
$ws.Range("L96").Value = 30997.334
$ws.Range("N96").Value = -36489.334
$ws.Range("H102").Value = 48561
$ws.Range("J102").Value = 48561
$ws.Range("L102").Value = 48561
$ws.Range("N102").Value = -55051
$ws.Range("H109").Value = 35273
$ws.Range("J109").Value = 35273
$ws.Range("L109").Value = 35273
$ws.Range("N109").Value = -38047
$ws.Range("H123").Value = 32875.332
$ws.Range("J123").Value = 32875.332
$ws.Range("L123").Value = 32875.332
$ws.Range("N123").Value = -42675.332
$ws.Range("H124").Value = 37598.4
$ws.Range("J124").Value = 37598.4
$ws.Range("L124").Value = 37598.4
$ws.Range("N124").Value = -47418.4
$ws.Range("H125").Value = 44846
$ws.Range("J125").Value = 44846
$ws.Range("L125").Value = 44846
$ws.Range("N125").Value = -54686
$ws.Range("H127").Value = 49558
$ws.Range("J127").Value = 49558
$ws.Range("L127").Value = 49558
$ws.Range("N127").Value = -59478
$ws.Range("H128").Value = 44429
$ws.Range("J128").Value = 44429
$ws.Range("L128").Value = 44429
$ws.Range("N128").Value = -54389
$ws.Range("H133").Value = 34549.5
$ws.Range("J133").Value = 34549.5
$ws.Range("L133").Value = 34549.5
$ws.Range("N133").Value = -39609.5
$ws.Range("H137").Value = 42000
$ws.Range("J137").Value = 42000
$ws.Range("L137").Value = 42000
$ws.Range("N137").Value = -52200

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H87").Value = 38700
$ws.Range("J87").Value = 38700
$ws.Range("L87").Value = 38700
$ws.Range("N87").Value = -41196
$ws.Range("H90").Value = 38700
$ws.Range("J90").Value = 38700
$ws.Range("L90").Value = 116100
$ws.Range("N90").Value = -128580
$ws.Range("H92").Value = 24478.166
$ws.Range("J92").Value = 24478.166
$ws.Range("L92").Value = 24478.166
$ws.Range("N92").Value = -29470.166
$ws.Range("H102").Value = 36168.5
$ws.Range("J102").Value = 36168.5
$ws.Range("L102").Value = 36168.5
$ws.Range("N102").Value = -42658.5
$ws.Range("H123").Value = 25679.25
$ws.Range("J123").Value = 24831.908
$ws.Range("L123").Value = 24831.908
$ws.Range("N123").Value = -34631.908
$ws.Range("H125").Value = 36254.5
$ws.Range("J125").Value = 36254.5
$ws.Range("L125").Value = 36254.5
$ws.Range("N125").Value = -46094.5
$ws.Range("H130").Value = 29286
$ws.Range("J130").Value = 29286
$ws.Range("L130").Value = 29286
$ws.Range("N130").Value = -39326
$ws.Range("H131").Value = 42263
$ws.Range("J131").Value = 42263
$ws.Range("L131").Value = 42263
$ws.Range("N131").Value = -52343
$ws.Range("H139").Value = 58500
$ws.Range("J139").Value = 58500
$ws.Range("L139").Value = 58500
$ws.Range("N139").Value = -68780
